$d = $word.ActiveDocument

# Replace the whole body content with the updated set of paragraphs. Every
# paragraph mark in the target revision is a plain, freshly authored one
# (no leftover proofing / language run-formatting), so we rebuild the body
# from scratch via WordprocessingML fragments rather than editing the
# existing (differently-formatted) paragraph marks in place.
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-Para([string]$runsXml) {
    return "<w:p " + $wNs + ">" + $runsXml + "</w:p>"
}

function New-Run([string]$text, [bool]$preserve) {
    $space = ""
    if ($preserve) { $space = ' xml:space="preserve"' }
    $escaped = $text -replace "&", "&amp;" -replace "<", "&lt;" -replace ">", "&gt;"
    return "<w:r " + $wNs + "><w:t" + $space + ">" + $escaped + "</w:t></w:r>"
}

$bodyXml = ""
$bodyXml += New-Para (New-Run "Tese" $false)
$bodyXml += New-Para ""
$bodyXml += New-Para (New-Run "Ssdfsdfsdf" $false)
$bodyXml += New-Para (New-Run "Sdfdsfdsfsdfs" $false)
$bodyXml += New-Para (New-Run "Dsfdsfdsdsfsdf" $false)
$bodyXml += New-Para (New-Run "Sdfsdfsdfsdf" $false)
$bodyXml += New-Para ""
$bodyXml += New-Para ((New-Run "Fazendo uma pequena alteração." $false) + (New-Run " " $true))
$bodyXml += New-Para ""
$bodyXml += New-Para ""

$full = $d.Range(0, $d.Content.End)
$full.InsertXML($bodyXml)

$d.Save()
